$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.985.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.69%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.429.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.52"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.78%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +7.33%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.56%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +22.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.49"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000220"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +72.68%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.974.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.23"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.66%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.481.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +13.97%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.977.10"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.58"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +24.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "89.99"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.07%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "32.59"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.70%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.55%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.76"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.53%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.62"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.64%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.39%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.90"

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "43.78"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.28%  "

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.171"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.17%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0501"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.68"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.64%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.74%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.12"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.78%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.07"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.85%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.43%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.86"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.118.72"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.131"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +16.88%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "PEPE"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₃0472"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +61.51%  "
